# Edit script: add 'diversity vs population size' sheet (sheet3),
# update selections on Sheet1 / Included species, and add the
# _xlnm._FilterDatabase defined name scoped to the new sheet.

$wb = $excel.ActiveWorkbook

# --- Update Sheet1 selection (drop old tabSelected/topLeftCell/activeCell) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1:A1048576").Select()

# --- Update 'Included species' selection ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1:A1048576").Select()

# --- Add the new worksheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "diversity vs population size"

# --- Column widths ---
$ws3.Columns.Item(1).ColumnWidth = 38.140625
$ws3.Columns.Item(2).ColumnWidth = 35.140625

# --- Header row ---
$ws3.Cells.Item(1, 1).Value = "Species with genetic diversity estimates"
$ws3.Cells.Item(1, 2).Value = "Species with population size estimates"

# --- Data rows: row|species|styleId|hasPopSizeCol|rowHeight ---
$data = @'
2|Aegilops tauschii|3|0|
3|Amaranthus hypochondriacus|3|0|
4|Ananas comosus|3|0|
5|Arabidopsis halleri|3|1|
6|Arabidopsis lyrata|3|1|
7|Arabidopsis suecica|3|1|
8|Arabidopsis thaliana|3|1|
9|Arabis alpina|3|0|
10|Arabis nemorensis|3|0|
11|Arachis duranensis|3|0|
12|Arachis hypogaea|3|0|
13|Arachis ipaensis|3|0|
14|Benincasa hispida|3|0|
15|Beta vulgaris|3|0|
16|Boechera stricta|3|1|
17|Brachypodium distachyon|3|0|
18|Brassica juncea|3|0|
19|Brassica napus|3|0|
20|Brassica oleracea capitata|3|0|
21|Brassica rapa|3|0|
22|Buddleja alternifolia|3|0|
23|Cajanus cajan|3|0|
24|Camellia sinensis|3|0|
25|Cannabis sativa|3|0|
26|Capsella grandiflora|3|1|
27|Capsella rubella|3|1|
28|Capsicum annuum|3|0|
29|Castanea mollissima|3|0|
30|Chenopodium quinoa|3|0|
31|Cicer arietinum |3|0|
32|Citrullus lanatus|3|0|
33|Coffea arabica|3|0|
34|Coffea canephora|3|0|
35|Cucumis melo|3|0|
36|Cucumis sativus|3|0|
37|Cucurbita pepo|3|0|
38|Digitaria exilis|3|0|
39|Dioscorea rotundata|3|0|
40|Elaeis guineensis|3|0|
41|Eleusine coracana|3|0|
42|Ficus carica|3|0|
43|Glycine max|3|0|
44|Glycine soja|3|1|
45|Gossypium arboreum|3|0|
46|Gossypium barbadense|3|0|
47|Gossypium hirsutum |3|0|
48|Hordeum vulgare|3|0|
49|Juglans regia|3|0|
50|Lactuca sativa|3|0|
51|Linum usitatissimum|3|0|
52|Lotus japonicus|3|0|
53|Lupinus angustifolius|3|0|
54|Macadamia integrifolia|3|0|
55|Malus domestica|3|0|
56|Malus sylvestris|3|0|
57|Mangifera indica|3|0|
58|Manihot esculenta|3|0|
59|Medicago truncatula|3|0|
60|Mimulus guttatus|3|1|
61|Momordica charantia|3|0|
62|Musa acuminata|3|0|
63|Nelumbo nucifera|3|0|
64|Olea europaea|3|0|
65|Oryza barthii|3|0|
66|Oryza brachyantha|5|0|
67|Oryza glaberrima|3|0|16.5
68|Oryza glumipatula|5|0|
69|Oryza longistaminata|5|0|
70|Oryza meridionalis|5|0|
71|Oryza punctata|5|0|
72|Oryza rufipogon|3|0|
73|Oryza sativa|3|0|
74|Panicum hallii|3|0|
75|Panicum virgatum|3|0|
76|Phaseolus vulgaris|3|0|
77|Phoenix dactylifera|3|0|
78|Pisum sativum|3|0|
79|Populus deltoides|3|0|
80|Populus trichocarpa|3|0|
81|Prunus armeniaca|3|0|
82|Prunus avium|3|0|
83|Prunus persica|3|0|
84|Quercus lobata|3|0|
85|Quercus robur|3|0|
86|Rhododendron griersonianum|3|0|
87|Salix dunnii|3|0|
88|Secale cereale|3|0|
89|Sesamum indicum|3|0|
90|Setaria italica |3|0|
91|Setaria viridis|3|0|
92|Solanum lycopersicum|3|0|
93|Solanum stenotomum|3|0|
94|Sorghum bicolor|3|0|
95|Spinacia oleracea |3|0|
96|Spirodela polyrhiza|3|0|
97|Striga hermonthica|3|0|
98|Tetracentron sinense|3|0|
99|Thlaspi arvense|3|0|
100|Triticum aestivum|3|0|
101|Triticum turgidum|3|0|
102|Vigna radiata|3|0|
103|Vigna umbellata|3|0|
104|Vigna unguiculata|3|0|
105|Vitis vinifera|3|0|
106|Xanthoceras sorbifolium|3|0|
107|Zea mays|3|0|
108|Ziziphus jujuba|3|0|
'@

$greenColor = 5296274
$yellowColor = 65535

$lines = $data -split "`n"
foreach ($line in $lines) {
    $parts = $line.Split("|")
    $r = [int]$parts[0]
    $species = $parts[1]
    $styleId = $parts[2]
    $hasB = $parts[3]
    $ht = $parts[4]

    $cellA = $ws3.Cells.Item($r, 1)
    $cellA.Value = $species
    if ($styleId -eq "3") {
        $cellA.Interior.Color = $greenColor
    } elseif ($styleId -eq "5") {
        $cellA.Interior.Color = $yellowColor
    }

    if ($hasB -eq "1") {
        $ws3.Cells.Item($r, 2).Value = $species
    }

    if ($ht -ne "") {
        $ws3.Rows.Item($r).RowHeight = [double]$ht
    }
}

# --- Hidden AutoFilter-database defined name scoped to the new sheet ---
$filterName = $ws3.Names.Add("_xlnm._FilterDatabase", "='diversity vs population size'!`$A`$1:`$B`$318")
$filterName.Visible = $False

# --- Final selection / active sheet ---
$ws3.Range("B44").Select()
$ws3.Activate()

